$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 25.603423448951439
$ws.Range("C2").Value = 16.190736754717673
$ws.Range("D2").Value = 0.63236608912862968
$ws.Range("E2").Value = 21.52706424336936
$ws.Range("F2").Value = 15.10978739599811
$ws.Range("G2").Value = 0.70189725943016767
$ws.Range("H2").Value = 257.5
$ws.Range("I2").Value = 226.5
$ws.Range("B3").Value = 26.656642092203942
$ws.Range("C3").Value = 16.073301208224848
$ws.Range("D3").Value = 0.60297546677590275
$ws.Range("E3").Value = 22.53843260815113
$ws.Range("F3").Value = 15.094749165498875
$ws.Range("G3").Value = 0.66973375779644007
$ws.Range("H3").Value = 269.5
$ws.Range("I3").Value = 237.5
$ws.Range("B4").Value = 27.803243624608427
$ws.Range("C4").Value = 16.159892658566807
$ws.Range("D4").Value = 0.58122328735284023
$ws.Range("E4").Value = 23.596090828205494
$ws.Range("F4").Value = 15.152150412448231
$ws.Range("G4").Value = 0.64214663872780864
$ws.Range("H4").Value = 280
$ws.Range("I4").Value = 247.5
$ws.Range("B5").Value = 28.583728265052901
$ws.Range("C5").Value = 15.938820537440163
$ws.Range("D5").Value = 0.55761866995241893
$ws.Range("E5").Value = 24.679274840690113
$ws.Range("F5").Value = 15.234799655303451
$ws.Range("G5").Value = 0.61731147911141127
$ws.Range("H5").Value = 289.5
$ws.Range("I5").Value = 257
$ws.Range("B6").Value = 29.83978691937951
$ws.Range("C6").Value = 16.209667250874858
$ws.Range("D6").Value = 0.54322329092596355
$ws.Range("E6").Value = 25.455337708048816
$ws.Range("F6").Value = 15.032760238492855
$ws.Range("G6").Value = 0.59055434309714905
$ws.Range("H6").Value = 299
$ws.Range("I6").Value = 265
$ws.Range("B7").Value = 30.976874103525958
$ws.Range("C7").Value = 16.263276750099465
$ws.Range("D7").Value = 0.52501348895782507
$ws.Range("E7").Value = 26.634861517978571
$ws.Range("F7").Value = 15.211843582363457
$ws.Range("G7").Value = 0.57112531154312329
$ws.Range("H7").Value = 307.5
$ws.Range("I7").Value = 273.5
$ws.Range("B8").Value = 31.679985721596857
$ws.Range("C8").Value = 15.931014898434375
$ws.Range("D8").Value = 0.50287317167488166
$ws.Range("E8").Value = 27.830410112127439
$ws.Range("F8").Value = 15.404591278285251
$ws.Range("G8").Value = 0.55351650285500154
$ws.Range("H8").Value = 316
$ws.Range("I8").Value = 281.5
$ws.Range("B9").Value = 32.893817835827981
$ws.Range("C9").Value = 16.102153727909759
$ws.Range("D9").Value = 0.48951914941205993
$ws.Range("E9").Value = 28.735609657137196
$ws.Range("F9").Value = 15.323526806979519
$ws.Range("G9").Value = 0.53325915092160026
$ws.Range("H9").Value = 324
$ws.Range("I9").Value = 288
$ws.Range("B10").Value = 34.141784424067168
$ws.Range("C10").Value = 16.306848117772411
$ws.Range("D10").Value = 0.477621436396787
$ws.Range("E10").Value = 29.644728774931963
$ws.Range("F10").Value = 15.284091722235457
$ws.Range("G10").Value = 0.51557536040471108
$ws.Range("H10").Value = 331.5
$ws.Range("I10").Value = 295
$ws.Range("B11").Value = 34.902943652677102
$ws.Range("C11").Value = 16.086617324795494
$ws.Range("D11").Value = 0.46089571942341401
$ws.Range("E11").Value = 30.945583760443384
$ws.Range("F11").Value = 15.53188159928167
$ws.Range("G11").Value = 0.50190947178496959
$ws.Range("H11").Value = 338.5
$ws.Range("I11").Value = 302